# Auto-generated Excel COM-interop edit script
# Applies cell-value corrections to the Leve profit tracking sheets
# (currentAveragePrice / Leve price / profit columns), per the scheduled
# market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = $null
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = $null
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("H132").Value = 3300.2273
$ws.Range("I132").Value = 2716.1765
$ws.Range("J132").Value = 5286
$ws.Range("K132").Value = 8148.529500000001
$ws.Range("L132").Value = 15858
$ws.Range("M132").Value = -5618.529500000001
$ws.Range("N132").Value = -20918
$ws.Range("H134").Value = 99100
$ws.Range("J134").Value = 99100
$ws.Range("L134").Value = 99100
$ws.Range("N134").Value = -109240
$ws.Range("H137").Value = 2905.4443
$ws.Range("I137").Value = 2487.5
$ws.Range("K137").Value = 7462.5
$ws.Range("M137").Value = -4912.5
$ws.Range("H138").Value = 3674.9285
$ws.Range("I138").Value = 2486.6667
$ws.Range("J138").Value = 3999
$ws.Range("K138").Value = 7460.000100000001
$ws.Range("L138").Value = 11997
$ws.Range("M138").Value = -2320.000100000001
$ws.Range("N138").Value = -22277
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2939
$ws.Range("I97").Value = 1723.75
$ws.Range("K97").Value = 1723.75
$ws.Range("M97").Value = -1227.75
$ws.Range("H110").Value = 995.6667
$ws.Range("I110").Value = 993.5
$ws.Range("K110").Value = 993.5
$ws.Range("M110").Value = 1051.5
$ws.Range("H139").Value = 80357.5
$ws.Range("I139").Value = 66000
$ws.Range("J139").Value = 94715
$ws.Range("K139").Value = 66000
$ws.Range("L139").Value = 94715
$ws.Range("M139").Value = -60860
$ws.Range("N139").Value = -104995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4999.5
$ws.Range("I20").Value = 4999.5
$ws.Range("K20").Value = 4999.5
$ws.Range("M20").Value = -4752.5
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = $null
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = $null
$ws.Range("H94").Value = 7202
$ws.Range("I94").Value = 6336
$ws.Range("J94").Value = 9800
$ws.Range("K94").Value = 6336
$ws.Range("L94").Value = 9800
$ws.Range("M94").Value = -5885
$ws.Range("N94").Value = -10702
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = $null
$ws.Range("H134").Value = 5727.25
$ws.Range("I134").Value = 4663.6
$ws.Range("K134").Value = 13990.8
$ws.Range("M134").Value = -11455.8
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6099.875
$ws.Range("I80").Value = 4799
$ws.Range("J80").Value = 6285.7144
$ws.Range("K80").Value = 14397
$ws.Range("L80").Value = 18857.1432
$ws.Range("M80").Value = -13461
$ws.Range("N80").Value = -20729.1432
$ws.Range("H83").Value = 6099.875
$ws.Range("I83").Value = 4799
$ws.Range("J83").Value = 6285.7144
$ws.Range("K83").Value = 43191
$ws.Range("L83").Value = 56571.4296
$ws.Range("M83").Value = -38511
$ws.Range("N83").Value = -65931.4296
$ws.Range("H109").Value = 1161.6666
$ws.Range("I109").Value = 1161.6666
$ws.Range("K109").Value = 3484.9998
$ws.Range("M109").Value = -2444.9998
$ws.Range("H116").Value = 3949.5
$ws.Range("I116").Value = 999
$ws.Range("J116").Value = 6900
$ws.Range("K116").Value = 2997
$ws.Range("L116").Value = 20700
$ws.Range("M116").Value = 445
$ws.Range("N116").Value = -27584
$ws.Range("H124").Value = 4987.5
$ws.Range("I124").Value = 4987.5
$ws.Range("K124").Value = 14962.5
$ws.Range("M124").Value = -10052.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3999.5
$ws.Range("I70").Value = 2499
$ws.Range("K70").Value = 2499
$ws.Range("M70").Value = -2229
$ws.Range("H73").Value = 3999.5
$ws.Range("I73").Value = 2499
$ws.Range("K73").Value = 2499
$ws.Range("M73").Value = -1563
$ws.Range("H95").Value = 62502
$ws.Range("J95").Value = 62502
$ws.Range("L95").Value = 62502
$ws.Range("N95").Value = -67994
$ws.Range("H97").Value = 1383.875
$ws.Range("I97").Value = 1799.2
$ws.Range("K97").Value = 1799.2
$ws.Range("M97").Value = -1303.2
$ws.Range("H107").Value = 607.5
$ws.Range("I107").Value = 212.5
$ws.Range("J107").Value = 1002.5
$ws.Range("K107").Value = 212.5
$ws.Range("L107").Value = 1002.5
$ws.Range("M107").Value = 1707.5
$ws.Range("N107").Value = -4842.5
$ws.Range("H110").Value = 99995
$ws.Range("J110").Value = 99995
$ws.Range("L110").Value = 99995
$ws.Range("N110").Value = -108175
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("H122").Value = 2526.1428
$ws.Range("J122").Value = 5749.5
$ws.Range("L122").Value = 17248.5
$ws.Range("N122").Value = -22148.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 2266
$ws.Range("I12").Value = 2266
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 2266
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -2096
$ws.Range("N12").Value = $null
$ws.Range("H112").Value = 99700
$ws.Range("J112").Value = 99700
$ws.Range("L112").Value = 99700
$ws.Range("N112").Value = -102654
$ws.Range("H122").Value = 4051.1428
$ws.Range("I122").Value = 4473.8
$ws.Range("J122").Value = 2994.5
$ws.Range("K122").Value = 13421.4
$ws.Range("L122").Value = 8983.5
$ws.Range("M122").Value = -10971.4
$ws.Range("N122").Value = -13883.5
$ws.Range("H132").Value = 3660.3333
$ws.Range("I132").Value = 3392.4
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 10177.2
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -7647.200000000001
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 50000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 50000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 50000
$ws.Range("M74").Value = $null
$ws.Range("N74").Value = -51872
$ws.Range("H77").Value = 50000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 50000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 150000
$ws.Range("M77").Value = $null
$ws.Range("N77").Value = -159360
